$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: new entry for github
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "github"
$ws.Range("D4").Value = "iskim0706"

# Row 5: new entry for gitlab
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "gitlab"
$ws.Range("D5").Value = "iskim0706"

# Match the final selection state shown in the diff
$ws.Range("D5").Select()
